$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 44: brand new experiment row (IA-PUCP / 5verbs - retrained / e-6 ...)
# ---------------------------------------------------------------------------
$ws.Range("A44").Value = "IA-PUCP"
$ws.Range("B44").Value = "5verbs - retrained"
$ws.Range("C44").Value = 58
$ws.Range("D44").Value = 10
$ws.Range("E44").Value = "e-6"
$ws.Range("F44").Value = 100
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 8
$ws.Range("I44").Value = 8
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1.518
$ws.Range("L44").Value = 0.3846
$ws.Range("M44").Value = 1.573
$ws.Range("N44").Value = 0.2609
$ws.Range("O44").Value = 0.62686567164179097

# A44 continues the yellow highlight used by A40:A43 above it
$ws.Range("A44").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Row 45: fill the previously blank placeholder row with the next experiment
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = "IA-PUCP"
$ws.Range("B45").Value = "5verbs - retrained"
$ws.Range("C45").Value = 59
$ws.Range("D45").Value = 10
$ws.Range("E45").Value = "e-6"
$ws.Range("F45").Value = 100
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 8
$ws.Range("I45").Value = 8
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1.511
$ws.Range("L45").Value = 0.3846
$ws.Range("M45").Value = 1.302
$ws.Range("N45").Value = 0.4783
$ws.Range("O45").Value = 0.63432835820895495

# A45 continues the yellow highlight as well
$ws.Range("A45").Interior.Color = 65535
# E45 already carried its (no-op) pre-existing style, nothing to change.
# N45 previously had the underline-font highlight style; that highlight is
# removed now that the cell holds real data again.
$ws.Range("N45").Font.Underline = $False

# ---------------------------------------------------------------------------
# View state: re-selected after the new rows were entered
# ---------------------------------------------------------------------------
$ws.Range("K49").Select()
